$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark exercises 3 to 6 (rows 3-6, column C) as done (TRUE)
$ws.Range("C3:C6").Value = $true

# Update the active selection to C6, matching the saved view state
$ws.Range("C6").Select()
